$wb = $excel.ActiveWorkbook

# --- "TODO CMS" sheet (sheet2.xml): add new TODO row for MailController strings ---
# (written first so its shared-string text takes the lower new sst index, matching
#  the authoring order captured in the target workbook)
$cms = $wb.Worksheets.Item("TODO CMS")
$cms.Cells.Item(16, 1).Value = "Strings (z.B. in MailController) in Datei auslagern"
$cms.Cells.Item(16, 2).Value = "offen"
$cms.Cells.Item(16, 2).Style = "Schlecht"

# --- "TODO" sheet (sheet1.xml): add new TODO row for MailController strings ---
$todo = $wb.Worksheets.Item("TODO")
$todo.Cells.Item(26, 1).Value = "Strings in MailController in Datei auslagern"
$todo.Cells.Item(26, 2).Value = "offen"
$todo.Cells.Item(26, 2).Style = "Schlecht"

# --- Update selection on "TODO CMS" (no longer the active tab) ---
$cms.Range("A29").Select()

# --- Make "TODO" the active sheet and update its selection ---
$todo.Select()
$todo.Range("C25").Select()
